# Fruta / hortaliza, semanal
# The data rows (2-22) get their D,L,M,N,O,P,Q,S,T values permuted across rows.
# Mapping: new row R receives the OLD values that used to live on row Map[R].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    2  = 21
    3  = 8
    4  = 22
    5  = 18
    6  = 12
    7  = 4
    8  = 19
    9  = 17
    10 = 5
    11 = 20
    12 = 13
    13 = 10
    14 = 7
    15 = 14
    16 = 2
    17 = 6
    18 = 3
    19 = 16
    20 = 9
    21 = 15
    22 = 11
}

# 1) Snapshot every source row's values BEFORE any writes happen, so the
#    permutation (a single 21-cycle) doesn't clobber data we still need.
$snapshot = @{}
foreach ($r in 2..22) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

# 2) Write the permuted values back out.
foreach ($r in 2..22) {
    $src = $snapshot[$map[$r]]
    $ws.Cells.Item($r, 4).Value = $src.D
    $ws.Cells.Item($r, 12).Value = $src.L
    $ws.Cells.Item($r, 13).Value = $src.M
    $ws.Cells.Item($r, 14).Value = $src.N
    $ws.Cells.Item($r, 15).Value = $src.O
    $ws.Cells.Item($r, 16).Value = $src.P
    $ws.Cells.Item($r, 17).Value = $src.Q
    $ws.Cells.Item($r, 19).Value = $src.S
    $ws.Cells.Item($r, 20).Value = $src.T
}
